# Add the new "Uldis test 19" project row to the sheet.
# (Sets A6/B6, which pushes the used range to A1:C6, grows the shared
# string table with the new text, and moves the active selection to the
# newly entered cell - mirroring what Excel does on manual data entry.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 916
$ws.Range("B6").Value = "Uldis test 19"

$ws.Range("B6").Select()
